$d = $word.ActiveDocument

# --- Step 1: insert "(720 soat onlayn va 144 soat oflayn shaklda) " before "pullik" ---
$found = $d.Content.Duplicate
$found.Find.Execute("pullik")
$ins = $found.Duplicate
$ins.Collapse(1)

$ins.InsertBefore("(720")
$ins.Font.Italic = $true
$ins.LanguageID = "uz-Cyrl-UZ"
$ins.LanguageIDFarEast = "x-none"
$ins.Collapse(0)

$ins.InsertBefore(" soat")
$ins.Font.Italic = $true
$ins.LanguageID = "en-US"
$ins.LanguageIDFarEast = "x-none"
$ins.Collapse(0)

$ins.InsertBefore(" ")
$ins.Font.Italic = $true
$ins.LanguageID = "uz-Cyrl-UZ"
$ins.LanguageIDFarEast = "x-none"
$ins.Collapse(0)

$ins.InsertBefore("onlaun va 144 soat oflayn shaklda")
$ins.Font.Italic = $true
$ins.LanguageID = "en-US"
$ins.LanguageIDFarEast = "x-none"
$ins.Collapse(0)

$ins.InsertBefore(") ")
$ins.Font.Italic = $true
$ins.LanguageID = "uz-Cyrl-UZ"
$ins.LanguageIDFarEast = "x-none"
$ins.Collapse(0)

Write-Output "step1 done"
